# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove header cell formatting (bold font + border + center/top alignment) ---
# and clear the stray "Unnamed: 0" header text in A1.
$ws.Range("A1").Value = ""
$ws.Range("A1:U1").ClearFormats()

# --- Row 3: Revisit count ---
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = ""
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 9
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 21

# --- Row 4: Fixation count ---
$ws.Range("C4").Value = 34
$ws.Range("D4").Value = ""
$ws.Range("I4").Value = 3
$ws.Range("K4").Value = 12
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 74

# --- Row 5: Dwell time (ms) ---
$ws.Range("C5").Value = 19487.37
$ws.Range("D5").Value = ""
$ws.Range("I5").Value = 2894.86
$ws.Range("K5").Value = 6473.78
$ws.Range("L5").Value = 3194.99
$ws.Range("M5").Value = 44873.12

# --- Row 6: Dwell time (%) ---
$ws.Range("C6").Value = 19.53
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = 0.62
$ws.Range("I6").Value = 2.9
$ws.Range("J6").Value = 0.87
$ws.Range("K6").Value = 6.49
$ws.Range("L6").Value = 3.2
$ws.Range("M6").Value = 44.96
$ws.Range("O6").Value = 0.25

# --- Row 7: Fixation duration (ms) ---
$ws.Range("C7").Value = 573.16
$ws.Range("D7").Value = ""
$ws.Range("I7").Value = 964.95
$ws.Range("K7").Value = 539.48
$ws.Range("L7").Value = 798.75
$ws.Range("M7").Value = 606.39

# --- Row 8: First fixation duration (ms) ---
$ws.Range("D8").Value = ""

# --- Remove now-unused trailing blank rows 10 and 11 ---
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
